$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the "Avance" (C) and "Fecha" (D) columns for rows 64..78 (Video 62..76):
# these items were previously at 0% with no completion date; they are now
# marked 100% complete with the date the work finished.

$updates = @(
    @{ Row = 64; Date = "10/12/2020" },
    @{ Row = 65; Date = "10/13/2020" },
    @{ Row = 66; Date = "10/13/2020" },
    @{ Row = 67; Date = "10/13/2020" },
    @{ Row = 68; Date = "10/14/2020" },
    @{ Row = 69; Date = "10/14/2020" },
    @{ Row = 70; Date = "10/14/2020" },
    @{ Row = 71; Date = "10/17/2020" },
    @{ Row = 72; Date = "10/17/2020" },
    @{ Row = 73; Date = "10/19/2020" },
    @{ Row = 74; Date = "10/19/2020" },
    @{ Row = 75; Date = "10/19/2020" },
    @{ Row = 76; Date = "10/20/2020" },
    @{ Row = 77; Date = "10/20/2020" },
    @{ Row = 78; Date = "10/20/2020" }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 3).Value = 100
    $ws.Cells.Item($r, 4).Value = $u.Date
}

# Apply the same date number format already used in the "Fecha" column (copy
# formatting from an existing date cell, e.g. D63, so the style is reused).
$ws.Range("D63").Copy() | Out-Null
$ws.Range("D64:D78").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Update the view: scroll position and active selection cell
$ws.Range("G84").Select()
$excel.ActiveWindow.ScrollRow = 72

$wb.Save()
